$d = $word.ActiveDocument

# 1. Replace the placeholder ID text in the first paragraph's first run.
[void]$d.Content.Find.Execute("**ID__AFFARS_pgi_5315_topic_44__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5315_406_3__ID**", 2)

$p1 = $d.Paragraphs.Item(1)

# 2. Remove the trailing space that used to live in its own run right
#    before the paragraph mark (that second run is being dropped entirely).
$r = $p1.Range
$spaceRange = $d.Range($r.Start + $r.Text.Length - 2, $r.Start + $r.Text.Length - 1)
$spaceRange.Delete()

# 3. Give the paragraph a (currently invisible, zero-width) border with
#    5pt spacing on every side, and bump its left indent from 120 twips
#    (6pt) to 225 twips (11.25pt) - matching the other paragraphs below it.
$pf = $p1.Range.ParagraphFormat
$pf.LeftIndent = 11.25

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
